$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B
$ws.Range("B2").Value = 2.436460462383138
$ws.Range("B3").Value = 2.281399856684232
$ws.Range("B4").Value = 2.187766516090051
$ws.Range("B5").Value = 2.150001498692745
$ws.Range("B6").Value = 2.143754149253198
$ws.Range("B7").Value = 2.187255625377077
$ws.Range("B8").Value = 2.3826659375747
$ws.Range("B9").Value = 2.778590196949324
$ws.Range("B10").Value = 3.077591619891223
$ws.Range("B11").Value = 3.215458578213543
$ws.Range("B12").Value = 3.267937295299305
$ws.Range("B13").Value = 3.256622900792081
$ws.Range("B14").Value = 3.219770560723362
$ws.Range("B15").Value = 3.197232960359941
$ws.Range("B16").Value = 3.068619378412791
$ws.Range("B17").Value = 2.99019674186917
$ws.Range("B18").Value = 2.945263881656729
$ws.Range("B19").Value = 2.930080066180039
$ws.Range("B20").Value = 2.998526941823911
$ws.Range("B21").Value = 3.230587573540163
$ws.Range("B22").Value = 3.383837728130061
$ws.Range("B23").Value = 3.301898343264611
$ws.Range("B24").Value = 2.994760381161996
$ws.Range("B25").Value = 2.670083293172581

# Column C
$ws.Range("C2").Value = 0.6052186298113611
$ws.Range("C3").Value = 0.5582328985321965
$ws.Range("C4").Value = 0.5297624186591747
$ws.Range("C5").Value = 0.5182542134313053
$ws.Range("C6").Value = 0.516348902231357
$ws.Range("C7").Value = 0.5296068374218521
$ws.Range("C8").Value = 0.5889383468946789
$ws.Range("C9").Value = 0.7083709899500263
$ws.Range("C10").Value = 0.798115261756152
$ws.Range("C11").Value = 0.8394023571956382
$ws.Range("C12").Value = 0.855105147556003
$ws.Range("C13").Value = 0.8517202050368269
$ws.Range("C14").Value = 0.840692856461942
$ws.Range("C15").Value = 0.8339472320639629
$ws.Range("C16").Value = 0.7954265047163176
$ws.Range("C17").Value = 0.7719148991418479
$ws.Range("C18").Value = 0.7584350167034586
$ws.Range("C19").Value = 0.7538783540692293
$ws.Range("C20").Value = 0.7744132478023289
$ws.Range("C21").Value = 0.8439299885851597
$ws.Range("C22").Value = 0.889761969656945
$ws.Range("C23").Value = 0.865263474557139
$ws.Range("C24").Value = 0.7732836284762925
$ws.Range("C25").Value = 0.675719189877384

# Column D
$ws.Range("D2").Value = 0.09191220859909066
$ws.Range("D3").Value = 0.09112208160775737
$ws.Range("D4").Value = 0.09067301728543242
$ws.Range("D5").Value = 0.09049912833443585
$ws.Range("D6").Value = 0.09047080542489283
$ws.Range("D7").Value = 0.09067063522897456
$ws.Range("D8").Value = 0.09163230631124719
$ws.Range("D9").Value = 0.09380301011465519
$ws.Range("D10").Value = 0.09556996235865256
$ws.Range("D11").Value = 0.09641088547384413
$ws.Range("D12").Value = 0.09673463204713784
$ws.Range("D13").Value = 0.09666467188254302
$ws.Range("D14").Value = 0.09643741411250772
$ws.Range("D15").Value = 0.09629890240277206
$ws.Range("D16").Value = 0.09551575042753768
$ws.Range("D17").Value = 0.09504479847460345
$ws.Range("D18").Value = 0.09477741582859522
$ws.Range("D19").Value = 0.09468748616166067
$ws.Range("D20").Value = 0.09509457047617076
$ws.Range("D21").Value = 0.09650402143771686
$ws.Range("D22").Value = 0.09745609983265524
$ws.Range("D23").Value = 0.09694513861781218
$ws.Range("D24").Value = 0.09507205804392527
$ws.Range("D25").Value = 0.09318544164072051

# Column F
$ws.Range("F2").Value = 3.27355071698986
$ws.Range("F3").Value = 3.212964882621577
$ws.Range("F4").Value = 3.177927133130737
$ws.Range("F5").Value = 3.164188257304943
$ws.Range("F6").Value = 3.16193937159693
$ws.Range("F7").Value = 3.177739668043444
$ws.Range("F8").Value = 3.252209257872181
$ws.Range("F9").Value = 3.415612518499046
$ws.Range("F10").Value = 3.546567403204421
$ws.Range("F11").Value = 3.608579818583905
$ws.Range("F12").Value = 3.632418756875495
$ws.Range("F13").Value = 3.627268688791958
$ws.Range("F14").Value = 3.610533893463241
$ws.Range("F15").Value = 3.600329891449263
$ws.Range("F16").Value = 3.542564299851193
$ws.Range("F17").Value = 3.507755729082049
$ws.Range("F18").Value = 3.487964139344768
$ws.Range("F19").Value = 3.481302284727718
$ws.Range("F20").Value = 3.511437383204225
$ws.Range("F21").Value = 3.615439599382711
$ws.Range("F22").Value = 3.685489680016502
$ws.Range("F23").Value = 3.647910674550246
$ws.Range("F24").Value = 3.509772221520336
$ws.Range("F25").Value = 3.369516051372756

# Column G
$ws.Range("G2").Value = 0.002528812161925748
$ws.Range("G3").Value = 0.00253549854366192
$ws.Range("G4").Value = 0.002539814162155079
$ws.Range("G5").Value = 0.002541625861157913
$ws.Range("G6").Value = 0.002541929902394132
$ws.Range("G7").Value = 0.002539838380266168
$ws.Range("G8").Value = 0.002531074132512459
$ws.Range("G9").Value = 0.002515545576644135
$ws.Range("G10").Value = 0.002505134419069322
$ws.Range("G11").Value = 0.002500611925394934
$ws.Range("G12").Value = 0.002498929871288504
$ws.Range("G13").Value = 0.002499290777774253
$ws.Range("G14").Value = 0.002500472931221041
$ws.Range("G15").Value = 0.002501201002879018
$ws.Range("G16").Value = 0.00250543425575101
$ws.Range("G17").Value = 0.002508085784542274
$ws.Range("G18").Value = 0.00250963099038418
$ws.Range("G19").Value = 0.002510157631298093
$ws.Range("G20").Value = 0.002507801444308755
$ws.Range("G21").Value = 0.002500124877575987
$ws.Range("G22").Value = 0.002495285576750712
$ws.Range("G23").Value = 0.002497852201408425
$ws.Range("G24").Value = 0.002507929929795871
$ws.Range("G25").Value = 0.002519570311447196

# Column J
$ws.Range("J2").Value = 0.3541570842887722
$ws.Range("J3").Value = 0.3435391397508027
$ws.Range("J4").Value = 0.3372583403832863
$ws.Range("J5").Value = 0.3347583493649751
$ws.Range("J6").Value = 0.3343468053085843
$ws.Range("J7").Value = 0.3372243844641503
$ws.Range("J8").Value = 0.3504461586090599
$ws.Range("J9").Value = 0.3782935567761996
$ws.Range("J10").Value = 0.3999620702061009
$ws.Range("J11").Value = 0.4100909037934173
$ws.Range("J12").Value = 0.4139661796176455
$ws.Range("J13").Value = 0.4131297934397224
$ws.Range("J14").Value = 0.4104089256116481
$ws.Range("J15").Value = 0.4087475054972032
$ws.Range("J16").Value = 0.3993056486146855
$ws.Range("J17").Value = 0.3935834151303226
$ws.Range("J18").Value = 0.3903176744259298
$ws.Range("J19").Value = 0.3892163176078611
$ws.Range("J20").Value = 0.3941899091626908
$ws.Range("J21").Value = 0.4112070277076469
$ws.Range("J22").Value = 0.4225604419021778
$ws.Range("J23").Value = 0.4164794934153093
$ws.Range("J24").Value = 0.3939156383554803
$ws.Range("J25").Value = 0.3705506311211622

# Column N
$ws.Range("N2").Value = 1.717516942719243
$ws.Range("N3").Value = 1.733059106012099
$ws.Range("N4").Value = 1.743258699992182
$ws.Range("N5").Value = 1.747579384048315
$ws.Range("N6").Value = 1.748306726905028
$ws.Range("N7").Value = 1.743316306250613
$ws.Range("N8").Value = 1.722739080689728
$ws.Range("N9").Value = 1.687635508305675
$ws.Range("N10").Value = 1.66509686897129
$ws.Range("N11").Value = 1.655561125820412
$ws.Range("N12").Value = 1.652054317493622
$ws.Range("N13").Value = 1.6528049216504
$ws.Range("N14").Value = 1.655270524843004
$ws.Range("N15").Value = 1.656794375211874
$ws.Range("N16").Value = 1.665734572417804
$ws.Range("N17").Value = 1.671403518444592
$ws.Range("N18").Value = 1.67473158906347
$ws.Range("N19").Value = 1.675869968581722
$ws.Range("N20").Value = 1.670793061926062
$ws.Range("N21").Value = 1.654543481887018
$ws.Range("N22").Value = 1.644531123185587
$ws.Range("N23").Value = 1.649818956665158
$ws.Range("N24").Value = 1.671068834732154
$ws.Range("N25").Value = 1.696564726526113
